$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 37, shifting existing rows 37-56 down to 38-57
$ws.Rows("37:37").Insert()

# Populate the new row 37 with the new data record
$ws.Range("A37").Value = 3
$ws.Range("B37").Value = "Femacal de La Calera"
$ws.Range("C37").Value = "Coquimbo"
$ws.Range("D37").Value = 44634
$ws.Range("D37").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E37").Value = 5
$ws.Range("F37").Value = 100112022
$ws.Range("G37").Value = "Arveja Verde"
$ws.Range("H37").Value = "Perfection"
$ws.Range("I37").Value = "Primera"
$ws.Range("J37").Value = 38
$ws.Range("K37").Value = 25000
$ws.Range("L37").Value = 25000
$ws.Range("M37").Value = 25000
$ws.Range("N37").Value = "`$/malla 25 kilos"
$ws.Range("O37").Value = "Provincia de Talca"
$ws.Range("P37").Value = 1000
$ws.Range("Q37").Value = 25
$ws.Range("R37").Value = "Hortaliza"
